# "new system time control and modification schemes"
# Update the recorded timestamp in C2 (date serial 45678 -> 45690, i.e.
# 2025-01-21 -> 2025-02-02) and leave the active selection parked on the
# cell that was just edited (matches the saved sheet selection moving
# from E2 to C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45690

$ws.Range("C2").Select()
